# Populate "visorInformacionTecnicaRed" test-case sheet: header row + two
# test-case rows, with the borders/fonts/alignment/column-row sizing that
# Excel produced when the rows were authored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
  "ID Caso",
  "Nombre/ Descripcion",
  "Tipo de prueba",
  "Módulo",
  "Precondiciones",
  "Pasos a seguir",
  "Datos de prueba",
  "Resultado esperado",
  "Resultado obtenido",
  "Estado (OK/FALLA)",
  "Automatizado",
  "Observaciones"
)

$row2 = @(
  "CP_INFTECRED_001",
  "ingreso a la vista visor de Informacion tecnica de red",
  "Positivo",
  "eCenter",
  "El usuario debe tener permisos para acceder a la vista",
  "1. Clic en módulo eCenter`n2. Scroll en el contenedor de aplicaciones`n3. Clic en ""visor de Informacion tecnica de red""",
  "N/A",
  "El sistema debe redirigido correctamente la vista visor de Informacion tecnica de red",
  "La vista se cargó sin errores",
  "OK",
  "SI",
  "N/A"
)

$row3 = @(
  "CP_INFTECRED_002",
  "", "", "", "", "", "", "", "", "", "", ""
)

# ---- values -------------------------------------------------------------
for ($c = 1; $c -le 12; $c++) {
  $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}
for ($c = 1; $c -le 12; $c++) {
  $ws.Cells.Item(2, $c).Value = $row2[$c - 1]
}
$ws.Cells.Item(3, 1).Value = $row3[0]

# ---- font helper ----------------------------------------------------------
function Set-AptosFont($rng, $bold) {
  $f = $rng.Font
  $f.Name = "Aptos Narrow"
  $f.Bold = $bold
  $f.Color = 0
}

# ---- header row (row 1): bold Aptos, full thin border, centered, wrap ---
$headerRng = $ws.Range("A1:L1")
Set-AptosFont $headerRng $true
$headerRng.Borders.LineStyle = 1
$headerRng.Borders.Weight = 2
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4108
$headerRng.WrapText = $true

# ---- column A (rows 2-3): Aptos regular, full thin border, vert-center wrap
$colARng = $ws.Range("A2:A3")
Set-AptosFont $colARng $false
$colARng.Borders.LineStyle = 1
$colARng.Borders.Weight = 2
$colARng.VerticalAlignment = -4108
$colARng.WrapText = $true

# ---- column F (rows 2-3): default font, full thin border, vert-center wrap
$colFRng = $ws.Range("F2:F3")
$colFRng.Font.Name = "Calibri"
$colFRng.Borders.LineStyle = 1
$colFRng.Borders.Weight = 2
$colFRng.VerticalAlignment = -4108
$colFRng.WrapText = $true

# ---- column C (rows 2-3): Aptos regular, border w/o left edge (black), vert-center wrap
$colCRng = $ws.Range("C2:C3")
Set-AptosFont $colCRng $false
$colCRng.Borders.Item(10).LineStyle = 1   # xlEdgeRight
$colCRng.Borders.Item(10).Weight = 2
$colCRng.Borders.Item(10).Color = 0
$colCRng.Borders.Item(8).LineStyle = 1    # xlEdgeTop
$colCRng.Borders.Item(8).Weight = 2
$colCRng.Borders.Item(8).Color = 0
$colCRng.Borders.Item(9).LineStyle = 1    # xlEdgeBottom
$colCRng.Borders.Item(9).Weight = 2
$colCRng.Borders.Item(9).Color = 0
$colCRng.VerticalAlignment = -4108
$colCRng.WrapText = $true

# ---- remaining data columns (B,D,E,G,H,I,J,K,L rows 2-3): Aptos regular,
#      full thin black border, vert-center wrap
$otherCols = @(2,4,5,7,8,9,10,11,12)
foreach ($c in $otherCols) {
  $rng = $ws.Range($ws.Cells.Item(2, $c), $ws.Cells.Item(3, $c))
  Set-AptosFont $rng $false
  $rng.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
  $rng.Borders.Item(7).Weight = 2
  $rng.Borders.Item(7).Color = 0
  $rng.Borders.Item(10).LineStyle = 1  # xlEdgeRight
  $rng.Borders.Item(10).Weight = 2
  $rng.Borders.Item(10).Color = 0
  $rng.Borders.Item(8).LineStyle = 1   # xlEdgeTop
  $rng.Borders.Item(8).Weight = 2
  $rng.Borders.Item(8).Color = 0
  $rng.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
  $rng.Borders.Item(9).Weight = 2
  $rng.Borders.Item(9).Color = 0
  $rng.VerticalAlignment = -4108
  $rng.WrapText = $true
}

# ---- column widths (character units ~= xml width - 0.8333333333333321) ---
$colWidths = @(
  20.166666666666668,
  20.022135416666668,
  17.307291666666668,
  15.592447916666668,
  17.022135416666668,
  18.592447916666668,
  22.307291666666668,
  23.592447916666668,
  23.166666666666668,
  26.592447916666668,
  29.736979166666668,
  33.87760416666667
)
for ($c = 1; $c -le 12; $c++) {
  $ws.Columns.Item($c).ColumnWidth = $colWidths[$c - 1]
}

# ---- row heights -----------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 60
$ws.Rows.Item(2).RowHeight = 120
$ws.Rows.Item(3).RowHeight = 84.75

# ---- view: active cell / selection / scroll -------------------------------
[void]$ws.Range("B3:L3").Select()
$excel.ActiveWindow.ScrollColumn = 2
